# Updates cryptos list with refreshed Price / Volume(1h) figures.
# For D-column values that look like plain numbers (e.g. "500.65", "1.00",
# "0.0000130") we prefix with a leading apostrophe so Excel stores them as
# text instead of silently normalising them to a Double (which would lose
# significant trailing zeros / exact decimal formatting). ClearFormats()
# afterwards strips the transient "quote-prefixed text" cell styling Excel
# applies so the cell keeps the workbook's original (unstyled) look.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.408.93'
$ws.Range("E2").Value = '  +0.83%  '
$ws.Range("D3").Value = '2.285.32'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = "'500.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("D6").Value = "'129.37"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.22%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'0.529"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = "'0.0956"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.50%  '
$ws.Range("D10").Value = "'0.152"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.22%  '
$ws.Range("D11").Value = "'0.334"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.72%  '
$ws.Range("D12").Value = "'4.73"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.41%  '
$ws.Range("D13").Value = '2.691.32'
$ws.Range("E13").Value = '  +0.56%  '
$ws.Range("D14").Value = "'23.10"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +7.34%  '
$ws.Range("D15").Value = '54.331.39'
$ws.Range("E15").Value = '  +0.58%  '
$ws.Range("D16").Value = "'0.0000130"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("D17").Value = '2.279.51'
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("D18").Value = "'10.28"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +4.79%  '
$ws.Range("D19").Value = "'4.13"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.48%  '
$ws.Range("D20").Value = "'305.51"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +2.24%  '
$ws.Range("D21").Value = "'6.42"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.33%  '
$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").Value = "'62.04"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.51%  '
$ws.Range("D24").Value = "'0.999"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.27%  '
$ws.Range("D25").Value = "'0.152"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.78%  '
$ws.Range("D26").Value = "'7.35"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.97%  '
$ws.Range("D27").Value = "'174.57"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +7.10%  '
$ws.Range("D28").Value = "'1.62"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.90%  '
$ws.Range("D29").Value = "'6.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.03%  '
$ws.Range("D30").Value = '0.0₃0691'
$ws.Range("E30").Value = '  +1.15%  '
$ws.Range("D31").Value = "'1.09"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.00%  '
$ws.Range("E32").Value = '  +0.09%  '
$ws.Range("D33").Value = "'17.83"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.10%  '
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").Value = "'0.936"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +9.47%  '
$ws.Range("E36").Value = '  +1.65%  '
$ws.Range("D37").Value = "'3.77"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +3.38%  '
$ws.Range("E38").Value = '  +1.54%  '
$ws.Range("E39").Value = '  -0.34%  '
$ws.Range("E40").Value = '  +1.24%  '
$ws.Range("D41").Value = "'3.40"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.87%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = "'4.90"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.93%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = "'125.22"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").Value = "'0.0496"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.47%  '
$ws.Range("D45").Value = "'0.0897"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.87%  '
$ws.Range("D46").Value = "'0.549"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.44%  '
$ws.Range("D47").Value = "'241.75"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.17%  '
$ws.Range("E48").Value = '  -0.40%  '
$ws.Range("D49").Value = "'0.0207"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.98%  '
$ws.Range("E50").Value = '  +1.26%  '
$ws.Range("D51").Value = "'16.42"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.17%  '
